$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - first sheet
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 579
$wsExpo.Range("F3").Value = 193
$wsExpo.Range("F4").Value = 400
$wsExpo.Range("I4").Value = "//i0.hdslb.com/bfs/openplatform/202406/CsYbpZmU1719311879090.jpeg"
$wsExpo.Range("F5").Value = 444
$wsExpo.Range("F7").Value = 2466
$wsExpo.Range("F8").Value = 425
$wsExpo.Range("F9").Value = 6490
$wsExpo.Range("F10").Value = 174
$wsExpo.Range("F11").Value = 420
$wsExpo.Range("F12").Value = 28

# Sheet "全部类型" (All types) - fourth sheet
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 579
$wsAll.Range("F3").Value = 193
$wsAll.Range("F4").Value = 400
$wsAll.Range("I4").Value = "//i0.hdslb.com/bfs/openplatform/202406/CsYbpZmU1719311879090.jpeg"
$wsAll.Range("F5").Value = 444
$wsAll.Range("F9").Value = 2466
$wsAll.Range("F10").Value = 425
$wsAll.Range("F11").Value = 6490
$wsAll.Range("F12").Value = 174
$wsAll.Range("F13").Value = 420
$wsAll.Range("F16").Value = 28
